# Weekly fruit/vegetable price update: insert 3 new daily price records
# (2022-01-17 data) for Ajo (garlic) at Femacal de La Calera, pushing the
# existing historical rows down by 3 positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows starting at row 314 (shifts old rows 314-338 down to 317-341)
$ws.Range("A314:R316").EntireRow.Insert()

# --- New row 314 ---
$ws.Range("A314").Value2 = 3
$ws.Range("B314").Value2 = "Femacal de La Calera"
$ws.Range("C314").Value2 = "Coquimbo"
$ws.Range("D314").Value2 = 44578
$ws.Range("E314").Value2 = 5
$ws.Range("F314").Value2 = 100112003
$ws.Range("G314").Value2 = "Ajo"
$ws.Range("H314").Value2 = "Chino"
$ws.Range("I314").Value2 = "1a (cosecha)"
$ws.Range("J314").Value2 = 110
$ws.Range("K314").Value2 = 16000
$ws.Range("L314").Value2 = 16500
$ws.Range("M314").Value2 = 16273
$ws.Range("N314").Value2 = "$/caja 10 kilos"
$ws.Range("O314").Value2 = "Llay Llay"
$ws.Range("P314").Value2 = 1627
$ws.Range("Q314").Value2 = 10
$ws.Range("R314").Value2 = "Hortaliza"

# --- New row 315 ---
$ws.Range("A315").Value2 = 3
$ws.Range("B315").Value2 = "Femacal de La Calera"
$ws.Range("C315").Value2 = "Coquimbo"
$ws.Range("D315").Value2 = 44578
$ws.Range("E315").Value2 = 5
$ws.Range("F315").Value2 = 100112003
$ws.Range("G315").Value2 = "Ajo"
$ws.Range("H315").Value2 = "Chino"
$ws.Range("I315").Value2 = "1a (cosecha)"
$ws.Range("J315").Value2 = 170
$ws.Range("K315").Value2 = 6500
$ws.Range("L315").Value2 = 7000
$ws.Range("M315").Value2 = 6765
$ws.Range("N315").Value2 = "$/trenza 50 unidades"
$ws.Range("O315").Value2 = "Llay Llay"
$ws.Range("P315").Value2 = 1353
$ws.Range("Q315").Value2 = 5
$ws.Range("R315").Value2 = "Hortaliza"

# --- New row 316 ---
$ws.Range("A316").Value2 = 3
$ws.Range("B316").Value2 = "Femacal de La Calera"
$ws.Range("C316").Value2 = "Coquimbo"
$ws.Range("D316").Value2 = 44578
$ws.Range("E316").Value2 = 5
$ws.Range("F316").Value2 = 100112003
$ws.Range("G316").Value2 = "Ajo"
$ws.Range("H316").Value2 = "Chino"
$ws.Range("I316").Value2 = "2a (cosecha)"
$ws.Range("J316").Value2 = 155
$ws.Range("K316").Value2 = 4000
$ws.Range("L316").Value2 = 4500
$ws.Range("M316").Value2 = 4242
$ws.Range("N316").Value2 = "$/trenza 50 unidades"
$ws.Range("O316").Value2 = "Llay Llay"
$ws.Range("P316").Value2 = 848
$ws.Range("Q316").Value2 = 5
$ws.Range("R316").Value2 = "Hortaliza"
